# Scheduled runner update: refresh market/profit figures across Sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) with latest pulled prices.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 1799.2858
$ws.Cells.Item(2, 9).Value = 100
$ws.Cells.Item(2, 10).Value = 2082.5
$ws.Cells.Item(2, 11).Value = 100
$ws.Cells.Item(2, 12).Value = 2082.5
$ws.Cells.Item(2, 13).Value = 13
$ws.Cells.Item(2, 14).Value = -2308.5

$ws.Cells.Item(6, 8).Value = 569.125
$ws.Cells.Item(6, 9).Value = 579
$ws.Cells.Item(6, 10).Value = 500
$ws.Cells.Item(6, 11).Value = 1737
$ws.Cells.Item(6, 12).Value = 1500
$ws.Cells.Item(6, 13).Value = -1625
$ws.Cells.Item(6, 14).Value = -1724

$ws.Cells.Item(17, 8).Value = 18333.334
$ws.Cells.Item(17, 10).Value = 18333.334
$ws.Cells.Item(17, 12).Value = 55000.00199999999
$ws.Cells.Item(17, 14).Value = -55336.00199999999

$ws.Cells.Item(33, 8).Value = 365.26086
$ws.Cells.Item(33, 9).Value = 304.6875
$ws.Cells.Item(33, 10).Value = 503.7143
$ws.Cells.Item(33, 11).Value = 304.6875
$ws.Cells.Item(33, 12).Value = 503.7143
$ws.Cells.Item(33, 13).Value = -75.6875
$ws.Cells.Item(33, 14).Value = -961.7143

$ws.Cells.Item(43, 8).Value = 7760.4
$ws.Cells.Item(43, 10).Value = 9625.5
$ws.Cells.Item(43, 12).Value = 9625.5
$ws.Cells.Item(43, 14).Value = -9763.5

$ws.Cells.Item(112, 8).Value = 2586.375
$ws.Cells.Item(112, 10).Value = 2586.375
$ws.Cells.Item(112, 12).Value = 7759.125
$ws.Cells.Item(112, 14).Value = -9975.125

$ws.Cells.Item(113, 8).Value = 7122.375
$ws.Cells.Item(113, 9).Value = 5250
$ws.Cells.Item(113, 11).Value = 5250
$ws.Cells.Item(113, 13).Value = -1996

$ws.Cells.Item(137, 8).Value = 11113806
$ws.Cells.Item(137, 9).Value = 31251812
$ws.Cells.Item(137, 10).Value = 3182.5344
$ws.Cells.Item(137, 11).Value = 93755436
$ws.Cells.Item(137, 12).Value = 9547.6032
$ws.Cells.Item(137, 13).Value = -93752886
$ws.Cells.Item(137, 14).Value = -14647.6032

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 4328.4814
$ws.Cells.Item(2, 9).Value = 537.5
$ws.Cells.Item(2, 10).Value = 11910.444
$ws.Cells.Item(2, 11).Value = 537.5
$ws.Cells.Item(2, 12).Value = 11910.444
$ws.Cells.Item(2, 13).Value = -424.5
$ws.Cells.Item(2, 14).Value = -12136.444

$ws.Cells.Item(61, 8).Value = 3092.551
$ws.Cells.Item(61, 9).Value = 2511.3262
$ws.Cells.Item(61, 11).Value = 2511.3262
$ws.Cells.Item(61, 13).Value = -2299.3262

$ws.Cells.Item(116, 8).Value = 4328.4814
$ws.Cells.Item(116, 9).Value = 537.5
$ws.Cells.Item(116, 10).Value = 11910.444
$ws.Cells.Item(116, 11).Value = 537.5
$ws.Cells.Item(116, 12).Value = 11910.444
$ws.Cells.Item(116, 13).Value = 1756.5
$ws.Cells.Item(116, 14).Value = -16498.444

$ws.Cells.Item(136, 8).Value = 3092.551
$ws.Cells.Item(136, 9).Value = 2511.3262
$ws.Cells.Item(136, 11).Value = 7533.9786
$ws.Cells.Item(136, 13).Value = -4983.9786

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 4328.4814
$ws.Cells.Item(3, 9).Value = 537.5
$ws.Cells.Item(3, 10).Value = 11910.444
$ws.Cells.Item(3, 11).Value = 537.5
$ws.Cells.Item(3, 12).Value = 11910.444
$ws.Cells.Item(3, 13).Value = -423.5
$ws.Cells.Item(3, 14).Value = -12138.444

$ws.Cells.Item(134, 8).Value = 2955
$ws.Cells.Item(134, 9).Value = 2285.0557
$ws.Cells.Item(134, 11).Value = 6855.1671
$ws.Cells.Item(134, 13).Value = -4320.1671

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 41135
$ws.Cells.Item(31, 9).Value = 1492.1177
$ws.Cells.Item(31, 11).Value = 1492.1177
$ws.Cells.Item(31, 13).Value = -1197.1177

$ws.Cells.Item(34, 8).Value = 41135
$ws.Cells.Item(34, 9).Value = 1492.1177
$ws.Cells.Item(34, 11).Value = 1492.1177
$ws.Cells.Item(34, 13).Value = -1290.1177

$ws.Cells.Item(41, 8).Value = 22764.75

$ws.Cells.Item(68, 8).Value = 80799.2
$ws.Cells.Item(68, 10).Value = 80799.2
$ws.Cells.Item(68, 12).Value = 80799.2
$ws.Cells.Item(68, 14).Value = -82297.2

$ws.Cells.Item(71, 8).Value = 80799.2
$ws.Cells.Item(71, 10).Value = 80799.2
$ws.Cells.Item(71, 12).Value = 242397.6
$ws.Cells.Item(71, 14).Value = -249885.6

$ws.Cells.Item(74, 8).Value = 80000
$ws.Cells.Item(74, 10).Value = 80000
$ws.Cells.Item(74, 12).Value = 80000
$ws.Cells.Item(74, 14).Value = -81748

$ws.Cells.Item(77, 8).Value = 80000
$ws.Cells.Item(77, 10).Value = 80000
$ws.Cells.Item(77, 12).Value = 240000
$ws.Cells.Item(77, 14).Value = -248736

$ws.Cells.Item(134, 8).Value = 3095.111
$ws.Cells.Item(134, 9).Value = 2133.5715
$ws.Cells.Item(134, 10).Value = 6460.5
$ws.Cells.Item(134, 11).Value = 6400.7145
$ws.Cells.Item(134, 12).Value = 19381.5
$ws.Cells.Item(134, 13).Value = -3865.7145
$ws.Cells.Item(134, 14).Value = -24451.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 114410.57
$ws.Cells.Item(2, 10).Value = 145603.19
$ws.Cells.Item(2, 12).Value = 873619.14
$ws.Cells.Item(2, 14).Value = -873845.14

$ws.Cells.Item(34, 8).Value = 5099.8
$ws.Cells.Item(34, 9).Value = 4266
$ws.Cells.Item(34, 10).Value = 5457.143
$ws.Cells.Item(34, 11).Value = 12798
$ws.Cells.Item(34, 12).Value = 16371.429
$ws.Cells.Item(34, 13).Value = -12714
$ws.Cells.Item(34, 14).Value = -16539.429

$ws.Cells.Item(122, 8).Value = 7661460
$ws.Cells.Item(122, 10).Value = 8406353
$ws.Cells.Item(122, 12).Value = 75657177
$ws.Cells.Item(122, 14).Value = -75662077

$ws.Cells.Item(124, 8).Value = 37041508
$ws.Cells.Item(124, 10).Value = 55560348
$ws.Cells.Item(124, 12).Value = 166681044
$ws.Cells.Item(124, 14).Value = -166690864

$ws.Cells.Item(131, 8).Value = 4841883.5
$ws.Cells.Item(131, 9).Value = 17858028
$ws.Cells.Item(131, 10).Value = 3297595
$ws.Cells.Item(131, 11).Value = 53574084
$ws.Cells.Item(131, 12).Value = 9892785
$ws.Cells.Item(131, 13).Value = -53569044
$ws.Cells.Item(131, 14).Value = -9902865

$ws.Cells.Item(133, 8).Value = 15879204
$ws.Cells.Item(133, 9).Value = 1197.375
$ws.Cells.Item(133, 11).Value = 3592.125
$ws.Cells.Item(133, 13).Value = 1467.875

$ws.Cells.Item(137, 8).Value = 52488.75
$ws.Cells.Item(137, 9).Value = 1157.7
$ws.Cells.Item(137, 10).Value = 103819.8
$ws.Cells.Item(137, 11).Value = 3473.1
$ws.Cells.Item(137, 12).Value = 311459.4
$ws.Cells.Item(137, 13).Value = 1626.9
$ws.Cells.Item(137, 14).Value = -321659.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(68, 8).Value = 50142.5
$ws.Cells.Item(68, 10).Value = 50295
$ws.Cells.Item(68, 12).Value = 50295
$ws.Cells.Item(68, 14).Value = -51917

$ws.Cells.Item(70, 8).Value = 10689.363
$ws.Cells.Item(70, 9).Value = 7630.8335
$ws.Cells.Item(70, 11).Value = 7630.8335
$ws.Cells.Item(70, 13).Value = -7360.8335

$ws.Cells.Item(71, 8).Value = 50142.5
$ws.Cells.Item(71, 10).Value = 50295
$ws.Cells.Item(71, 12).Value = 150885
$ws.Cells.Item(71, 14).Value = -158997

$ws.Cells.Item(73, 8).Value = 10689.363
$ws.Cells.Item(73, 9).Value = 7630.8335
$ws.Cells.Item(73, 11).Value = 7630.8335
$ws.Cells.Item(73, 13).Value = -6694.8335

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 10435.786
$ws.Cells.Item(7, 9).Value = 7044.5557
$ws.Cells.Item(7, 10).Value = 16540
$ws.Cells.Item(7, 11).Value = 7044.5557
$ws.Cells.Item(7, 12).Value = 16540
$ws.Cells.Item(7, 13).Value = -6932.5557
$ws.Cells.Item(7, 14).Value = -16764

$ws.Cells.Item(68, 8).Value = 7684
$ws.Cells.Item(68, 9).Value = 4249.5
$ws.Cells.Item(68, 10).Value = 8828.833000000001
$ws.Cells.Item(68, 11).Value = 4249.5
$ws.Cells.Item(68, 12).Value = 8828.833000000001
$ws.Cells.Item(68, 13).Value = -3500.5
$ws.Cells.Item(68, 14).Value = -10326.833

$ws.Cells.Item(71, 8).Value = 7684
$ws.Cells.Item(71, 9).Value = 4249.5
$ws.Cells.Item(71, 10).Value = 8828.833000000001
$ws.Cells.Item(71, 11).Value = 21247.5
$ws.Cells.Item(71, 12).Value = 44144.165
$ws.Cells.Item(71, 13).Value = -17503.5
$ws.Cells.Item(71, 14).Value = -51632.165

$ws.Cells.Item(122, 8).Value = 216466.38
$ws.Cells.Item(122, 9).Value = 271963.8
$ws.Cells.Item(122, 11).Value = 815891.3999999999
$ws.Cells.Item(122, 13).Value = -813441.3999999999

$ws.Cells.Item(126, 8).Value = 10435.786
$ws.Cells.Item(126, 9).Value = 7044.5557
$ws.Cells.Item(126, 10).Value = 16540
$ws.Cells.Item(126, 11).Value = 21133.6671
$ws.Cells.Item(126, 12).Value = 49620
$ws.Cells.Item(126, 13).Value = -18663.6671
$ws.Cells.Item(126, 14).Value = -54560

$ws.Cells.Item(134, 8).Value = 58409.668
$ws.Cells.Item(134, 10).Value = 58409.668
$ws.Cells.Item(134, 12).Value = 58409.668
$ws.Cells.Item(134, 14).Value = -68549.66800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 7999.222
$ws.Cells.Item(62, 9).Value = 7990
$ws.Cells.Item(62, 10).Value = 8000.375
$ws.Cells.Item(62, 11).Value = 7990
$ws.Cells.Item(62, 12).Value = 8000.375
$ws.Cells.Item(62, 13).Value = -7366
$ws.Cells.Item(62, 14).Value = -9248.375

$ws.Cells.Item(65, 8).Value = 7999.222
$ws.Cells.Item(65, 9).Value = 7990
$ws.Cells.Item(65, 10).Value = 8000.375
$ws.Cells.Item(65, 11).Value = 39950
$ws.Cells.Item(65, 12).Value = 40001.875
$ws.Cells.Item(65, 13).Value = -36830
$ws.Cells.Item(65, 14).Value = -46241.875

$ws.Cells.Item(81, 8).Value = 4339.8
$ws.Cells.Item(81, 9).Value = 2710.6667
$ws.Cells.Item(81, 10).Value = 19002
$ws.Cells.Item(81, 11).Value = 5421.3334
$ws.Cells.Item(81, 12).Value = 38004
$ws.Cells.Item(81, 13).Value = -4360.3334
$ws.Cells.Item(81, 14).Value = -40126

$ws.Cells.Item(84, 8).Value = 4339.8
$ws.Cells.Item(84, 9).Value = 2710.6667
$ws.Cells.Item(84, 10).Value = 19002
$ws.Cells.Item(84, 11).Value = 27106.667
$ws.Cells.Item(84, 12).Value = 190020
$ws.Cells.Item(84, 13).Value = -21802.667
$ws.Cells.Item(84, 14).Value = -200628

$ws.Cells.Item(132, 8).Value = 2867.182
$ws.Cells.Item(132, 9).Value = 2153.7
$ws.Cells.Item(132, 11).Value = 6461.099999999999
$ws.Cells.Item(132, 13).Value = -3931.099999999999
